$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 19608824
$ws.Range("I19").Value = 83333780
$ws.Range("K19").Value = 83333780
$ws.Range("M19").Value = -83333605

# Row 28
$ws.Range("H28").Value = 705.13635
$ws.Range("I28").Value = 459.82352
$ws.Range("J28").Value = 1539.2
$ws.Range("K28").Value = 459.82352
$ws.Range("L28").Value = 1539.2
$ws.Range("M28").Value = 25.17648000000003
$ws.Range("N28").Value = -2509.2

# Row 33
$ws.Range("H33").Value = 389.52173
$ws.Range("I33").Value = 380.75
$ws.Range("J33").Value = 399.0909
$ws.Range("K33").Value = 380.75
$ws.Range("L33").Value = 399.0909
$ws.Range("M33").Value = -151.75
$ws.Range("N33").Value = -857.0908999999999

# Row 55
$ws.Range("H55").Value = 80.111115
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 106
$ws.Range("H106").Value = 2116
$ws.Range("I106").Value = 1150
$ws.Range("J106").Value = 5980
$ws.Range("K106").Value = 1150
$ws.Range("L106").Value = 5980
$ws.Range("M106").Value = -519
$ws.Range("N106").Value = -7242

# Row 125
$ws.Range("H125").Value = 16772.572
$ws.Range("I125").Value = 9800
$ws.Range("J125").Value = 17934.666
$ws.Range("K125").Value = 88200
$ws.Range("L125").Value = 161411.994
$ws.Range("M125").Value = -85740
$ws.Range("N125").Value = -166331.994

# Row 137
$ws.Range("H137").Value = 3881.5144
$ws.Range("I137").Value = 1975.421
$ws.Range("K137").Value = 5926.263
$ws.Range("M137").Value = -3376.263

# Row 138
$ws.Range("H138").Value = 3094.99
$ws.Range("I138").Value = 1065.3704
$ws.Range("J138").Value = 3845.6711
$ws.Range("K138").Value = 3196.1112
$ws.Range("L138").Value = 11537.0133
$ws.Range("M138").Value = 1943.8888
$ws.Range("N138").Value = -21817.0133

# Row 141
$ws.Range("H141").Value = 2089.4
$ws.Range("J141").Value = 3273.75
$ws.Range("L141").Value = 9821.25
$ws.Range("N141").Value = -20181.25

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 9777.147000000001
$ws.Range("J61").Value = 20475.363
$ws.Range("L61").Value = 20475.363
$ws.Range("N61").Value = -20899.363

# Row 74
$ws.Range("H74").Value = 115870.16
$ws.Range("I74").Value = 133197.44
$ws.Range("K74").Value = 133197.44
$ws.Range("M74").Value = -132323.44

# Row 77
$ws.Range("H77").Value = 115870.16
$ws.Range("I77").Value = 133197.44
$ws.Range("K77").Value = 665987.2
$ws.Range("M77").Value = -661619.2

# Row 132
$ws.Range("H132").Value = 6781.9355
$ws.Range("I132").Value = 3045.375
$ws.Range("J132").Value = 8081.609
$ws.Range("K132").Value = 9136.125
$ws.Range("L132").Value = 24244.827
$ws.Range("M132").Value = -6606.125
$ws.Range("N132").Value = -29304.827

# Row 134
$ws.Range("H134").Value = 65430
$ws.Range("J134").Value = 65430
$ws.Range("L134").Value = 65430
$ws.Range("N134").Value = -75570

# Row 136
$ws.Range("H136").Value = 9777.147000000001
$ws.Range("J136").Value = 20475.363
$ws.Range("L136").Value = 61426.08900000001
$ws.Range("N136").Value = -66526.08900000001

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 5437.8276
$ws.Range("I105").Value = 5160.3477
$ws.Range("J105").Value = 6501.5
$ws.Range("K105").Value = 5160.3477
$ws.Range("L105").Value = 6501.5
$ws.Range("M105").Value = -3413.3477
$ws.Range("N105").Value = -9995.5

# Row 134
$ws.Range("H134").Value = 41083.92
$ws.Range("I134").Value = 2817.65
$ws.Range("J134").Value = 168638.17
$ws.Range("K134").Value = 8452.950000000001
$ws.Range("L134").Value = 505914.51
$ws.Range("M134").Value = -5917.950000000001
$ws.Range("N134").Value = -510984.51

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1667
$ws.Range("I16").Value = 1500.5
$ws.Range("K16").Value = 1500.5
$ws.Range("M16").Value = -1213.5

# Row 31
$ws.Range("H31").Value = 2975.2593
$ws.Range("I31").Value = 2595.182
$ws.Range("J31").Value = 3236.5625
$ws.Range("K31").Value = 2595.182
$ws.Range("L31").Value = 3236.5625
$ws.Range("M31").Value = -2300.182
$ws.Range("N31").Value = -3826.5625

# Row 34
$ws.Range("H34").Value = 2975.2593
$ws.Range("I34").Value = 2595.182
$ws.Range("J34").Value = 3236.5625
$ws.Range("K34").Value = 2595.182
$ws.Range("L34").Value = 3236.5625
$ws.Range("M34").Value = -2393.182
$ws.Range("N34").Value = -3640.5625

# Row 107
$ws.Range("H107").Value = 578.9231
$ws.Range("I107").Value = 515.375
$ws.Range("J107").Value = 680.6
$ws.Range("K107").Value = 515.375
$ws.Range("L107").Value = 680.6
$ws.Range("M107").Value = 1404.625
$ws.Range("N107").Value = -4520.6

# Row 113
$ws.Range("H113").Value = 1667
$ws.Range("I113").Value = 1500.5
$ws.Range("K113").Value = 1500.5
$ws.Range("M113").Value = 669.5

# Row 132
$ws.Range("H132").Value = 2901.85
$ws.Range("I132").Value = 2391.8125
$ws.Range("J132").Value = 4942
$ws.Range("K132").Value = 7175.4375
$ws.Range("L132").Value = 14826
$ws.Range("M132").Value = -4645.4375
$ws.Range("N132").Value = -19886

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1245.6364
$ws.Range("I16").Value = 1070.2
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 1070.2
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -900.2
$ws.Range("N16").Value = -3340

# Row 122
$ws.Range("H122").Value = 5456.4653
$ws.Range("I122").Value = 4566.7646
$ws.Range("J122").Value = 6716.875
$ws.Range("K122").Value = 13700.2938
$ws.Range("L122").Value = 20150.625
$ws.Range("M122").Value = -11250.2938
$ws.Range("N122").Value = -25050.625

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 58968.7
$ws.Range("J46").Value = 58968.7
$ws.Range("L46").Value = 58968.7
$ws.Range("N46").Value = -59430.7

# Row 132
$ws.Range("H132").Value = 3208.625
$ws.Range("I132").Value = 2930.7727
$ws.Range("J132").Value = 3548.2222
$ws.Range("K132").Value = 8792.3181
$ws.Range("L132").Value = 10644.6666
$ws.Range("M132").Value = -6262.3181
$ws.Range("N132").Value = -15704.6666

# Row 134
$ws.Range("H134").Value = 58968.7
$ws.Range("J134").Value = 58968.7
$ws.Range("L134").Value = 176906.1
$ws.Range("N134").Value = -181976.1

# Row 135
$ws.Range("H135").Value = 500025000
$ws.Range("J135").Value = 500025000
$ws.Range("L135").Value = 500025000
$ws.Range("N135").Value = -500035140
